# Adds a new "optimization_parameters" worksheet (with species name / taxon id
# rows, per the commit message) to the expression-sheet test workbook, and makes
# it the active/selected sheet — mirroring the author's edit.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet at the end of the tab strip -----------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "optimization_parameters"

# --- populate the optimization-parameters table -------------------------------
$ws.Cells.Item(1, 1).Value = "optimization_parameter"
$ws.Cells.Item(1, 2).Value = "value"

$ws.Cells.Item(2, 1).Value = "alpha"
$ws.Cells.Item(2, 2).Value = 0.002
$ws.Cells.Item(2, 2).NumberFormat = "0.00E+00"

$ws.Cells.Item(3, 1).Value = "kk_max"
$ws.Cells.Item(3, 2).Value = 1

$ws.Cells.Item(4, 1).Value = "MaxIter"
$ws.Cells.Item(4, 2).Value = 100000000
$ws.Cells.Item(4, 2).NumberFormat = "0.00E+00"

$ws.Cells.Item(5, 1).Value = "TolFun"
$ws.Cells.Item(5, 2).Value = 0.000001
$ws.Cells.Item(5, 2).NumberFormat = "0.00E+00"

$ws.Cells.Item(6, 1).Value = "MaxFunEval"
$ws.Cells.Item(6, 2).Value = 100000000
$ws.Cells.Item(6, 2).NumberFormat = "0.00E+00"

$ws.Cells.Item(7, 1).Value = "TolX"
$ws.Cells.Item(7, 2).Value = 0.000001
$ws.Cells.Item(7, 2).NumberFormat = "0.00E+00"

$ws.Cells.Item(8, 1).Value = "production_function"
$ws.Cells.Item(8, 2).Value = "Sigmoid"

$ws.Cells.Item(9, 1).Value = "L_curve"
$ws.Cells.Item(9, 2).Value = 0

$ws.Cells.Item(10, 1).Value = "estimate_params"
$ws.Cells.Item(10, 2).Value = 1

$ws.Cells.Item(11, 1).Value = "make_graphs"
$ws.Cells.Item(11, 2).Value = 1

$ws.Cells.Item(12, 1).Value = "fix_P"
$ws.Cells.Item(12, 2).Value = 0

$ws.Cells.Item(13, 1).Value = "fix_b"
$ws.Cells.Item(13, 2).Value = 0

$ws.Cells.Item(14, 1).Value = "expression_timepoints"
$ws.Cells.Item(14, 2).Value = 15
$ws.Cells.Item(14, 3).Value = 30
$ws.Cells.Item(14, 4).Value = 60

$ws.Cells.Item(15, 1).Value = "Strain"
$ws.Cells.Item(15, 2).Value = "wt"
$ws.Cells.Item(15, 3).Value = "dgln3"

$ws.Cells.Item(16, 1).Value = "simulation_timepoints"
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 10
$ws.Cells.Item(16, 5).Value = 15

$ws.Cells.Item(17, 1).Value = "species"
$ws.Cells.Item(17, 2).Value = "Saccharomyces cerevisiae"

$ws.Cells.Item(18, 1).Value = "taxon_id"
$ws.Cells.Item(18, 2).Value = 559292

# --- make the new sheet the active / selected tab, like the author's edit -----
$ws.Activate()

# --- iterative-calculation setting referenced by the workbook's calcPr --------
$excel.Iteration  = $true
$excel.MaxChange  = 0.0001
